$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '57.425.06'
$ws.Range("E2").Value = '  +0.35%  '
$ws.Range("D3").Value = '2.367.63'
$ws.Range("E3").Value = '  +0.13%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = "'518.53"
$ws.Range("E5").Value = '  -0.32%  '
$ws.Range("D6").Value = "'135.61"
$ws.Range("E6").Value = '  +0.44%  '
$ws.Range("D7").Value = "'0.996"
$ws.Range("E7").Value = '  -0.43%  '
$ws.Range("D8").Value = "'0.541"
$ws.Range("E8").Value = '  +0.16%  '
$ws.Range("E9").Value = '  -0.98%  '
$ws.Range("D10").Value = "'5.52"
$ws.Range("E10").Value = '  +5.24%  '
$ws.Range("E11").Value = '  -0.93%  '
$ws.Range("E12").Value = '  -0.02%  '
$ws.Range("D13").Value = "'24.38"
$ws.Range("E13").Value = '  +1.79%  '
$ws.Range("D14").Value = '2.792.33'
$ws.Range("E14").Value = '  +0.29%  '
$ws.Range("D15").Value = '57.446.80'
$ws.Range("E15").Value = '  +0.63%  '
$ws.Range("E16").Value = '  +0.12%  '
$ws.Range("D17").Value = '2.393.48'
$ws.Range("E17").Value = '  +1.05%  '
$ws.Range("D18").Value = "'10.60"
$ws.Range("E18").Value = '  +0.32%  '
$ws.Range("D19").Value = "'330.27"
$ws.Range("E19").Value = '  +2.19%  '
$ws.Range("E20").Value = '  -1.07%  '
$ws.Range("E21").Value = '  -0.12%  '
$ws.Range("D22").Value = "'0.998"
$ws.Range("E22").Value = '  -0.24%  '
$ws.Range("D23").Value = "'61.50"
$ws.Range("E23").Value = '  +0.02%  '
$ws.Range("D24").Value = "'8.96"
$ws.Range("E24").Value = '  +14.96%  '
$ws.Range("E25").Value = '  +3.56%  '
$ws.Range("D26").Value = "'0.996"
$ws.Range("E26").Value = '  +0.33%  '
$ws.Range("E27").Value = '  +10.61%  '
$ws.Range("D28").Value = '0.0₃0748'
$ws.Range("E28").Value = '  +0.56%  '
$ws.Range("E29").Value = '  +1.10%  '
$ws.Range("D30").Value = "'167.04"
$ws.Range("E30").Value = '  -2.70%  '
$ws.Range("E31").Value = '  -0.23%  '
$ws.Range("D32").Value = "'18.62"
$ws.Range("E32").Value = '  +1.12%  '
$ws.Range("E33").Value = '  -0.06%  '
$ws.Range("E34").Value = '  +3.36%  '
$ws.Range("D35").Value = "'0.994"
$ws.Range("E35").Value = '  -0.48%  '
$ws.Range("D36").Value = "'0.920"
$ws.Range("E36").Value = '  -3.86%  '
$ws.Range("E37").Value = '  +0.41%  '
$ws.Range("D38").Value = "'1.61"
$ws.Range("E38").Value = '  +6.16%  '
$ws.Range("D39").Value = "'38.89"
$ws.Range("E39").Value = '  +3.47%  '
$ws.Range("D40").Value = "'150.57"
$ws.Range("E40").Value = '  +7.15%  '
$ws.Range("D41").Value = "'0.388"
$ws.Range("E41").Value = '  +1.21%  '
$ws.Range("E42").Value = '  +1.81%  '
$ws.Range("D43").Value = "'289.89"
$ws.Range("E43").Value = '  +3.84%  '
$ws.Range("D44").Value = "'5.32"
$ws.Range("E44").Value = '  +3.02%  '
$ws.Range("E45").Value = '  +1.07%  '
$ws.Range("E46").Value = '  -0.48%  '
$ws.Range("E47").Value = '  +0.90%  '
$ws.Range("B48").Value = 'Polygon'
$ws.Range("C48").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D48").Value = "'0.390"
$ws.Range("E48").Value = '  +1.75%  '
$ws.Range("B49").Value = 'InjectiveProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D49").Value = "'18.22"
$ws.Range("E49").Value = '  +4.96%  '
$ws.Range("D50").Value = "'17.80"
$ws.Range("E50").Value = '  +4.45%  '
$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").Value = "'0.0219"
$ws.Range("E51").Value = '  +1.49%  '
